$wb = $excel.ActiveWorkbook

# --- RAF-capacity: correctly factor in capacity factor improvements / adjust RAF for capacity ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")

$wsCap.Range("B6").Value = 1      # hydro
$wsCap.Range("B7").Value = 0.3    # onshore wind
$wsCap.Range("B8").Value = 0.3    # solar pv
$wsCap.Range("B9").Value = 0.3    # solar thermal
$wsCap.Range("B10").Value = 0.3   # biomass
$wsCap.Range("B15").Value = 0.3   # offshore wind

# Clear the leftover formatting on the header row/column (row height, header cell style)
$wsCap.Rows.Item(1).AutoFit()
$wsCap.Range("A1").Style = "Normal"

# Make RAF-capacity the active sheet with the whole sheet selected (matches the saved view state)
$wsCap.Activate()
$wsCap.Cells.Select()

$wb.Save()
